$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B through AD) were swapped between each
# other. Column A (the running index) stays put.
$pairs = @(
    @(139, 140),
    @(175, 176),
    @(231, 232),
    @(267, 268)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $rng1.Value()
    $vals2 = $rng2.Value()

    $rng1.Value = $vals2
    $rng2.Value = $vals1
}
